$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Fill in the "SONUÇ" column (C2:C11) with the IF formula comparing distance to 500
$ws.Range("C2").Formula = '=IF(B2<500,"YAKIN","UZAK")'
$ws.Range("C3:C11").Formula = '=IF(B3<500,"YAKIN","UZAK")'

# Fill in the student info box (Numara / Ad Soyad / Bölüm)
$ws.Range("G11").Value = 20215070019
$ws.Range("G12").Value = "Kübra Çabuk"
$ws.Range("G13").Value = "YBS"

# Update the active selection to match the final state of the workbook
$ws.Range("G13:I13").Select()

$wb.Save()
